# COVID testing data source details - update for 2020-08-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Source URL string updates ---
$ws.Range("D4").Value = "https://www.health.gov.au/sites/default/files/documents/2020/08/coronavirus-covid-19-at-a-glance-20-august-2020.pdf"
$ws.Range("D17").Value = "https://www.facebook.com/Mshpci/posts/1673896616109359"
$ws.Range("D22").Value = "https://files.ssi.dk/Data-Epidemiologiske-Rapport-20082020-jh34"
$ws.Range("D24").Value = "https://www.facebook.com/nayibbukele/posts/3164190897000413"
$ws.Range("D26").Value = "https://twitter.com/EPHIEthiopia/status/1296485928240918531"
$ws.Range("Q26").Value = "https://twitter.com/EPHIEthiopia"
$ws.Range("D40").Value = "http://irangov.ir/detail/345517"
$ws.Range("D43").Value = "https://govextra.gov.il/media/24739/covid-19-data-israel-13082020.csv"
$ws.Range("D46").Value = "https://www.mhlw.go.jp/stf/newpage_13073.html"
$ws.Range("D47").Value = "https://www.mhlw.go.jp/content/10906000/000661150.pdf"
$ws.Range("D49").Value = "https://twitter.com/MOH_Kenya/status/1295717267376943106"
$ws.Range("D54").Value = "http://covid-19.moh.gov.my/terkini/082020/situasi-terkini-20-ogos-2020"
$ws.Range("D69").Value = "https://www.dge.gob.pe/portal/docs/tools/coronavirus/coronavirus180820.pdf"
$ws.Range("D71").Value = "https://twitter.com/MZ_GOV_PL/status/1296371398386913281"
$ws.Range("D72").Value = "https://twitter.com/MZ_GOV_PL/status/1296371398386913281"
$ws.Range("D76").Value = "https://rospotrebnadzor.ru/about/info/news/news_details.php?ELEMENT_ID=15210"
$ws.Range("D77").Value = "https://twitter.com/RwandaHealth/status/1296181973946425350"
$ws.Range("D92").Value = "https://ddc.moph.go.th/viralpneumonia/file/situation/situation-no230-200863.pdf"
$ws.Range("D93").Value = "https://ddc.moph.go.th/viralpneumonia/file/situation/situation-no230-200863.pdf"
$ws.Range("D97").Value = "https://twitter.com/MinofHealthUG/status/1296425945960767489/photo/2"
$ws.Range("D100").Value = "https://assets.publishing.service.gov.uk/government/uploads/system/uploads/attachment_data/file/910857/2020-08-20_COVID-19_UK_testing_time_series.csv"
$ws.Range("D103").Value = "https://www.gub.uy/ministerio-salud-publica/comunicacion/noticias/informe-situacion-sobre-coronavirus-covid-19-uruguay-20-agosto"
$ws.Range("D105").Value = "https://twitter.com/MoHCCZim/status/1296519577103466508"

# --- Numeric data updates ---
# Row 2
$ws.Range("C2").Value = 44063
$ws.Range("G2").Value = 177
$ws.Range("H2").Value = 901255
$ws.Range("I2").Value = 19.941
$ws.Range("J2").Value = 4391
$ws.Range("K2").Value = 0.097
$ws.Range("L2").Value = 11114
$ws.Range("M2").Value = 0.246
$ws.Range("N2").Value = 0.579
$ws.Range("O2").Value = 1.727

# Row 3
$ws.Range("C3").Value = 44063
$ws.Range("G3").Value = 177
$ws.Range("H3").Value = 1095233
$ws.Range("I3").Value = 24.233
$ws.Range("J3").Value = 5714
$ws.Range("K3").Value = 0.126
$ws.Range("L3").Value = 13231
$ws.Range("M3").Value = 0.293
$ws.Range("N3").Value = 0.486
$ws.Range("O3").Value = 2.056

# Row 4
$ws.Range("C4").Value = 44063
$ws.Range("G4").Value = 138
$ws.Range("H4").Value = 5508831
$ws.Range("I4").Value = 216.034
$ws.Range("J4").Value = 68336
$ws.Range("K4").Value = 2.68
$ws.Range("L4").Value = 57007
$ws.Range("O4").Value = 213.853

# Row 6
$ws.Range("C6").Value = 44064
$ws.Range("G6").Value = 141
$ws.Range("H6").Value = 1011805
$ws.Range("I6").Value = 594.626
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("L6").Value = 9572
$ws.Range("M6").Value = 5.625
$ws.Range("N6").ClearContents()
$ws.Range("O6").ClearContents()

# Row 7
$ws.Range("C7").Value = 44061
$ws.Range("G7").Value = 167
$ws.Range("H7").Value = 1378819
$ws.Range("I7").Value = 8.372
$ws.Range("J7").Value = 14630
$ws.Range("K7").Value = 0.089
$ws.Range("L7").Value = 12976
$ws.Range("M7").Value = 0.079
$ws.Range("N7").Value = 0.205
$ws.Range("O7").Value = 4.874

# Row 9
$ws.Range("C9").Value = 44062
$ws.Range("G9").Value = 172
$ws.Range("H9").Value = 2072738
$ws.Range("I9").Value = 178.844
$ws.Range("J9").Value = 18492
$ws.Range("K9").Value = 1.596
$ws.Range("L9").Value = 17796
$ws.Range("M9").Value = 1.536
$ws.Range("N9").Value = 0.023
$ws.Range("O9").Value = 44.363

# Row 12
$ws.Range("C12").Value = 44064
$ws.Range("G12").Value = 109
$ws.Range("H12").Value = 358723
$ws.Range("I12").Value = 51.626
$ws.Range("J12").Value = 6844
$ws.Range("K12").Value = 0.985
$ws.Range("L12").Value = 4839
$ws.Range("M12").Value = 0.696
$ws.Range("N12").ClearContents()
$ws.Range("O12").ClearContents()

# Row 13
$ws.Range("C13").Value = 44064
$ws.Range("G13").Value = 153
$ws.Range("H13").Value = 4974215
$ws.Range("I13").Value = 131.795
$ws.Range("J13").Value = 48769
$ws.Range("K13").Value = 1.292
$ws.Range("L13").Value = 48444
$ws.Range("M13").Value = 1.284
$ws.Range("N13").ClearContents()
$ws.Range("O13").ClearContents()

# Row 14
$ws.Range("C14").Value = 44063
$ws.Range("G14").Value = 142
$ws.Range("H14").Value = 2113632
$ws.Range("I14").Value = 110.568
$ws.Range("J14").Value = 26278
$ws.Range("K14").Value = 1.375
$ws.Range("L14").Value = 25863
$ws.Range("M14").Value = 1.353
$ws.Range("N14").Value = 0.066
$ws.Range("O14").Value = 15.253

# Row 15
$ws.Range("C15").Value = 44063
$ws.Range("G15").Value = 169
$ws.Range("H15").Value = 2309447
$ws.Range("I15").Value = 45.388
$ws.Range("J15").Value = 31932
$ws.Range("K15").Value = 0.628
$ws.Range("L15").Value = 35159
$ws.Range("M15").Value = 0.691
$ws.Range("N15").Value = 0.324
$ws.Range("O15").Value = 3.09

# Row 16
$ws.Range("C16").Value = 44062
$ws.Range("G16").Value = 162
$ws.Range("H16").Value = 109207
$ws.Range("I16").Value = 21.438
$ws.Range("J16").Value = 1535
$ws.Range("K16").Value = 0.301
$ws.Range("L16").Value = 1820
$ws.Range("M16").Value = 0.357

# Row 17
$ws.Range("C17").Value = 44062
$ws.Range("G17").Value = 126
$ws.Range("H17").Value = 116785
$ws.Range("I17").Value = 4.427
$ws.Range("J17").Value = 1447
$ws.Range("K17").Value = 0.055
$ws.Range("L17").Value = 882
$ws.Range("M17").Value = 0.033
$ws.Range("N17").Value = 0.049
$ws.Range("O17").Value = 20.376

# Row 18
$ws.Range("C18").Value = 44063
$ws.Range("G18").Value = 162
$ws.Range("H18").Value = 144521
$ws.Range("I18").Value = 35.204
$ws.Range("J18").Value = 2397
$ws.Range("K18").Value = 0.584
$ws.Range("L18").Value = 1575
$ws.Range("M18").Value = 0.384
$ws.Range("N18").Value = 0.109
$ws.Range("O18").Value = 9.157

# Row 19
$ws.Range("C19").Value = 44062
$ws.Range("G19").Value = 151
$ws.Range("H19").Value = 346439
$ws.Range("I19").Value = 30.586
$ws.Range("J19").Value = 5116
$ws.Range("K19").Value = 0.452
$ws.Range("L19").Value = 4563
$ws.Range("M19").Value = 0.403
$ws.Range("O19").Value = 101.4

# Row 20
$ws.Range("C20").Value = 44062
$ws.Range("G20").Value = 200
$ws.Range("H20").Value = 820151
$ws.Range("I20").Value = 76.585
$ws.Range("J20").Value = 7185
$ws.Range("K20").Value = 0.671
$ws.Range("L20").Value = 6237
$ws.Range("M20").Value = 0.582
$ws.Range("O20").Value = 25.682

# Row 22
$ws.Range("C22").Value = 44062
$ws.Range("G22").Value = 193
$ws.Range("H22").Value = 2035007
$ws.Range("I22").Value = 351.336
$ws.Range("J22").Value = 7704
$ws.Range("K22").Value = 1.33
$ws.Range("L22").Value = 27990
$ws.Range("M22").Value = 4.832
$ws.Range("O22").Value = 218.672

# Row 24
$ws.Range("C24").Value = 44062
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = 288307
$ws.Range("I24").Value = 44.449
$ws.Range("J24").Value = 2414
$ws.Range("K24").Value = 0.372
$ws.Range("L24").Value = 2464
$ws.Range("M24").Value = 0.38
$ws.Range("N24").Value = 0.127
$ws.Range("O24").Value = 7.865

# Row 26
$ws.Range("C26").Value = 44063
$ws.Range("G26").Value = 85
$ws.Range("H26").Value = 694093
$ws.Range("I26").Value = 6.038
$ws.Range("J26").Value = 21456
$ws.Range("K26").Value = 0.187
$ws.Range("L26").Value = 20568
$ws.Range("M26").Value = 0.179
$ws.Range("N26").Value = 0.062
$ws.Range("O26").Value = 16.105

# Row 29
$ws.Range("C29").Value = 44060
$ws.Range("G29").Value = 97
$ws.Range("J29").Value = 127008
$ws.Range("K29").Value = 1.946
$ws.Range("L29").Value = 82448
$ws.Range("M29").Value = 1.263
$ws.Range("N29").Value = 0.036
$ws.Range("O29").Value = 27.996

# Row 32
$ws.Range("C32").Value = 44061
$ws.Range("G32").Value = 117
$ws.Range("H32").Value = 431272
$ws.Range("I32").Value = 13.879
$ws.Range("J32").Value = 1499
$ws.Range("K32").Value = 0.048
$ws.Range("L32").Value = 1383
$ws.Range("M32").Value = 0.045
$ws.Range("N32").Value = 0.149
$ws.Range("O32").Value = 6.718

# Row 40
$ws.Range("C40").Value = 44063
$ws.Range("G40").Value = 120
$ws.Range("H40").Value = 2963741
$ws.Range("I40").Value = 35.286
$ws.Range("J40").Value = 23901
$ws.Range("K40").Value = 0.285
$ws.Range("L40").Value = 25126
$ws.Range("M40").Value = 0.299
$ws.Range("N40").Value = 0.08
$ws.Range("O40").Value = 12.442

# Row 42
$ws.Range("C42").Value = 44063
$ws.Range("G42").Value = 156
$ws.Range("H42").Value = 747521
$ws.Range("I42").Value = 151.388
$ws.Range("J42").Value = 12416
$ws.Range("K42").Value = 2.514
$ws.Range("L42").Value = 8546
$ws.Range("M42").Value = 1.731
$ws.Range("N42").Value = 0.012
$ws.Range("O42").Value = 84.375

# Row 43
$ws.Range("C43").Value = 44056
$ws.Range("G43").Value = 194
$ws.Range("H43").Value = 2050053
$ws.Range("I43").Value = 236.849
$ws.Range("J43").Value = 28904
$ws.Range("K43").Value = 3.339
$ws.Range("L43").Value = 21889
$ws.Range("M43").Value = 2.529
$ws.Range("O43").Value = 14.978

# Row 44
$ws.Range("C44").Value = 44063
$ws.Range("G44").Value = 124
$ws.Range("H44").Value = 4600949
$ws.Range("I44").Value = 76.097
$ws.Range("J44").Value = 49662
$ws.Range("K44").Value = 0.821
$ws.Range("L44").Value = 31185
$ws.Range("M44").Value = 0.516
$ws.Range("N44").Value = 0.016
$ws.Range("O44").Value = 61.233

# Row 45
$ws.Range("C45").Value = 44063
$ws.Range("G45").Value = 179
$ws.Range("H45").Value = 7790596
$ws.Range("I45").Value = 128.851
$ws.Range("J45").Value = 77442
$ws.Range("K45").Value = 1.281
$ws.Range("L45").Value = 52833
$ws.Range("M45").Value = 0.874
$ws.Range("O45").Value = 103.739

# Row 46
$ws.Range("C46").Value = 44063
$ws.Range("G46").Value = 189
$ws.Range("H46").Value = 1259422
$ws.Range("I46").Value = 9.958
$ws.Range("J46").Value = 23242
$ws.Range("K46").Value = 0.184
$ws.Range("L46").Value = 24871
$ws.Range("M46").Value = 0.197
$ws.Range("N46").Value = 0.046
$ws.Range("O46").Value = 21.608

# Row 47
$ws.Range("C47").Value = 44061
$ws.Range("G47").Value = 47
$ws.Range("H47").Value = 1547172
$ws.Range("I47").Value = 12.233
$ws.Range("J47").Value = 18957
$ws.Range("K47").Value = 0.15
$ws.Range("L47").Value = 20690
$ws.Range("M47").Value = 0.164
$ws.Range("N47").Value = 0.053
$ws.Range("O47").Value = 18.76

# Row 49
$ws.Range("C49").Value = 44061
$ws.Range("G49").Value = 135
$ws.Range("H49").Value = 398585
$ws.Range("I49").Value = 7.413
$ws.Range("J49").Value = 4019
$ws.Range("K49").Value = 0.075
$ws.Range("L49").Value = 5155
$ws.Range("N49").Value = 0.095
$ws.Range("O49").Value = 10.499

# Row 51
$ws.Range("C51").Value = 44064
$ws.Range("G51").Value = 175
$ws.Range("H51").Value = 234520
$ws.Range("I51").Value = 124.335
$ws.Range("J51").Value = 1602
$ws.Range("K51").Value = 0.849
$ws.Range("L51").Value = 1713
$ws.Range("M51").Value = 0.908
$ws.Range("N51").ClearContents()
$ws.Range("O51").ClearContents()

# Row 53
$ws.Range("C53").Value = 44062
$ws.Range("G53").Value = 175
$ws.Range("H53").Value = 466964
$ws.Range("I53").Value = 745.977
$ws.Range("J53").Value = 2781
$ws.Range("K53").Value = 4.443
$ws.Range("L53").Value = 1660
$ws.Range("M53").Value = 2.652
$ws.Range("N53").Value = 0.022
$ws.Range("O53").Value = 45.214

# Row 54
$ws.Range("C54").Value = 44063
$ws.Range("G54").Value = 163
$ws.Range("H54").Value = 1152140
$ws.Range("I54").Value = 35.597
$ws.Range("J54").Value = 10461
$ws.Range("K54").Value = 0.323
$ws.Range("L54").Value = 9586
$ws.Range("M54").Value = 0.296
$ws.Range("O54").Value = 554.562

# Row 57
$ws.Range("C57").Value = 44060
$ws.Range("G57").Value = 230
$ws.Range("H57").Value = 1138828
$ws.Range("I57").Value = 8.833
$ws.Range("J57").Value = 7608
$ws.Range("K57").Value = 0.059
$ws.Range("L57").Value = 9386
$ws.Range("M57").Value = 0.073
$ws.Range("N57").Value = 0.637
$ws.Range("O57").Value = 1.569

# Row 59
$ws.Range("C59").Value = 44061
$ws.Range("G59").Value = 133
$ws.Range("H59").Value = 139514
$ws.Range("I59").Value = 2.564
$ws.Range("J59").Value = 1765
$ws.Range("L59").Value = 1530
$ws.Range("M59").Value = 0.028
$ws.Range("O59").Value = 669.375

# Row 62
$ws.Range("C62").Value = 44063
$ws.Range("G62").Value = 165
$ws.Range("H62").Value = 673220
$ws.Range("I62").Value = 139.608
$ws.Range("J62").Value = 15714
$ws.Range("K62").Value = 3.259
$ws.Range("L62").Value = 21258
$ws.Range("N62").Value = 0
$ws.Range("O62").Value = 2254.636

# Row 66
$ws.Range("C66").Value = 44064
$ws.Range("G66").Value = 161
$ws.Range("H66").Value = 2389365
$ws.Range("I66").Value = 10.817
$ws.Range("J66").Value = 25613
$ws.Range("K66").Value = 0.116
$ws.Range("L66").Value = 22851
$ws.Range("M66").Value = 0.103
$ws.Range("N66").ClearContents()
$ws.Range("O66").ClearContents()

# Row 68
$ws.Range("C68").Value = 44062
$ws.Range("G68").Value = 164
$ws.Range("H68").Value = 162323
$ws.Range("I68").Value = 22.758
$ws.Range("J68").Value = 2657
$ws.Range("K68").Value = 0.373
$ws.Range("L68").Value = 2544
$ws.Range("M68").Value = 0.357
$ws.Range("N68").Value = 0.173
$ws.Range("O68").Value = 5.769

# Row 69
$ws.Range("C69").Value = 44061
$ws.Range("G69").Value = 167
$ws.Range("H69").Value = 506901
$ws.Range("I69").Value = 15.374
$ws.Range("J69").Value = 9044
$ws.Range("K69").Value = 0.274
$ws.Range("L69").Value = 8538
$ws.Range("M69").Value = 0.259

# Row 71
$ws.Range("C71").Value = 44063
$ws.Range("G71").Value = 114
$ws.Range("H71").Value = 2357772
$ws.Range("I71").Value = 62.298
$ws.Range("J71").Value = 24646
$ws.Range("K71").Value = 0.651
$ws.Range("L71").Value = 21510
$ws.Range("M71").Value = 0.568
$ws.Range("O71").Value = 30.511

# Row 72
$ws.Range("C72").Value = 44063
$ws.Range("G72").Value = 165
$ws.Range("H72").Value = 2463734
$ws.Range("I72").Value = 65.098
$ws.Range("J72").Value = 25599
$ws.Range("K72").Value = 0.676
$ws.Range("L72").Value = 22286
$ws.Range("M72").Value = 0.589
$ws.Range("N72").Value = 0.032
$ws.Range("O72").Value = 31.611

# Row 74
$ws.Range("C74").Value = 44063
$ws.Range("G74").Value = 153
$ws.Range("H74").Value = 572273
$ws.Range("I74").Value = 198.633
$ws.Range("J74").Value = 6260
$ws.Range("K74").Value = 2.173
$ws.Range("L74").Value = 4810
$ws.Range("M74").Value = 1.67
$ws.Range("N74").Value = 0.06
$ws.Range("O74").Value = 16.685

# Row 76
$ws.Range("C76").Value = 44063
$ws.Range("G76").Value = 166
$ws.Range("H76").Value = 33814105
$ws.Range("I76").Value = 231.707
$ws.Range("J76").Value = 304832
$ws.Range("K76").Value = 2.089
$ws.Range("L76").Value = 273007
$ws.Range("O76").Value = 55.201

# Row 77
$ws.Range("C77").Value = 44062
$ws.Range("G77").Value = 136
$ws.Range("H77").Value = 345920
$ws.Range("I77").Value = 26.707
$ws.Range("J77").Value = 3306
$ws.Range("K77").Value = 0.255
$ws.Range("L77").Value = 4483
$ws.Range("M77").Value = 0.346
$ws.Range("N77").Value = 0.013
$ws.Range("O77").Value = 77.293

# Row 78
$ws.Range("C78").Value = 44063
$ws.Range("G78").Value = 105
$ws.Range("H78").Value = 4501104
$ws.Range("I78").Value = 129.291
$ws.Range("J78").Value = 61620
$ws.Range("K78").Value = 1.77
$ws.Range("L78").Value = 61321
$ws.Range("M78").Value = 1.761
$ws.Range("N78").Value = 0.022
$ws.Range("O78").Value = 44.486

# Row 79
$ws.Range("C79").Value = 44063
$ws.Range("G79").Value = 171
$ws.Range("H79").Value = 138713
$ws.Range("I79").Value = 8.284
$ws.Range("J79").Value = 1519
$ws.Range("K79").Value = 0.091
$ws.Range("L79").Value = 1445
$ws.Range("M79").Value = 0.086
$ws.Range("N79").Value = 0.085
$ws.Range("O79").Value = 11.775

# Row 85
$ws.Range("C85").Value = 44063
$ws.Range("G85").Value = 173
$ws.Range("H85").Value = 3480283
$ws.Range("I85").Value = 58.681
$ws.Range("J85").Value = 24612
$ws.Range("K85").Value = 0.415
$ws.Range("L85").Value = 23541
$ws.Range("M85").Value = 0.397
$ws.Range("N85").Value = 0.165
$ws.Range("O85").Value = 6.072

# Row 91
$ws.Range("C91").Value = 44063
$ws.Range("G91").Value = 217
$ws.Range("H91").Value = 85563
$ws.Range("I91").Value = 3.593
$ws.Range("J91").Value = 214
$ws.Range("L91").Value = 188
$ws.Range("O91").Value = 263.2

# Row 92
$ws.Range("C92").Value = 44063
$ws.Range("G92").Value = 158
$ws.Range("H92").Value = 401680
$ws.Range("I92").Value = 5.755
$ws.Range("J92").Value = 2286
$ws.Range("K92").Value = 0.033
$ws.Range("L92").Value = 1506
$ws.Range("M92").Value = 0.022
$ws.Range("O92").Value = 319.455

# Row 93
$ws.Range("C93").Value = 44063
$ws.Range("G93").Value = 63
$ws.Range("H93").Value = 789951
$ws.Range("I93").Value = 11.317
$ws.Range("J93").Value = 2286
$ws.Range("K93").Value = 0.033
$ws.Range("L93").Value = 1506
$ws.Range("M93").Value = 0.022
$ws.Range("O93").Value = 319.455

# Row 94
$ws.Range("C94").Value = 44062
$ws.Range("G94").Value = 163
$ws.Range("H94").Value = 54709
$ws.Range("I94").Value = 6.608
$ws.Range("J94").Value = 862
$ws.Range("K94").Value = 0.104
$ws.Range("L94").Value = 698
$ws.Range("M94").Value = 0.084
$ws.Range("N94").Value = 0.021
$ws.Range("O94").Value = 47.437

# Row 95
$ws.Range("C95").Value = 44060
$ws.Range("G95").Value = 160
$ws.Range("H95").Value = 113758
$ws.Range("I95").Value = 9.625
$ws.Range("J95").Value = 1710
$ws.Range("K95").Value = 0.145
$ws.Range("L95").Value = 1365
$ws.Range("M95").Value = 0.115
$ws.Range("N95").Value = 0.043
$ws.Range("O95").Value = 23.305

# Row 96
$ws.Range("H96").Value = 6061930
$ws.Range("I96").Value = 71.876
$ws.Range("L96").Value = 77150
$ws.Range("M96").Value = 0.915
$ws.Range("N96").Value = 0.016
$ws.Range("O96").Value = 61.961

# Row 97
$ws.Range("C97").Value = 44062
$ws.Range("G97").Value = 46
$ws.Range("H97").Value = 336713
$ws.Range("I97").Value = 7.361
$ws.Range("J97").Value = 3046
$ws.Range("K97").Value = 0.067
$ws.Range("L97").Value = 4030
$ws.Range("M97").Value = 0.088
$ws.Range("N97").Value = 0.01
$ws.Range("O97").Value = 97.276

# Row 98
$ws.Range("C98").Value = 44064
$ws.Range("G98").Value = 119
$ws.Range("H98").Value = 1378519
$ws.Range("I98").Value = 31.521
$ws.Range("J98").Value = 21821
$ws.Range("K98").Value = 0.499
$ws.Range("L98").Value = 17954
$ws.Range("M98").Value = 0.411
$ws.Range("N98").ClearContents()
$ws.Range("O98").ClearContents()

# Row 99
$ws.Range("C99").Value = 44063
$ws.Range("G99").Value = 174
$ws.Range("H99").Value = 6265918
$ws.Range("I99").Value = 633.535
$ws.Range("J99").Value = 72283
$ws.Range("K99").Value = 7.308
$ws.Range("L99").Value = 71509
$ws.Range("M99").Value = 7.23
$ws.Range("O99").Value = 235.116

# Row 100
$ws.Range("C100").Value = 44062
$ws.Range("G100").Value = 142
$ws.Range("H100").Value = 12446595
$ws.Range("I100").Value = 183.346
$ws.Range("J100").Value = 175916
$ws.Range("K100").Value = 2.591
$ws.Range("L100").Value = 162256
$ws.Range("M100").Value = 2.39
$ws.Range("O100").Value = 151.5

# Row 102
$ws.Range("C102").Value = 44063
$ws.Range("G102").Value = 167
$ws.Range("H102").Value = 69580676
$ws.Range("I102").Value = 210.212
$ws.Range("J102").Value = 635809
$ws.Range("K102").Value = 1.921
$ws.Range("L102").Value = 709806
$ws.Range("M102").Value = 2.144
$ws.Range("N102").Value = 0.067
$ws.Range("O102").Value = 14.935

# Row 103
$ws.Range("C103").Value = 44063
$ws.Range("G103").Value = 130
$ws.Range("H103").Value = 154431
$ws.Range("I103").Value = 44.457
$ws.Range("J103").Value = 2310
$ws.Range("K103").Value = 0.665
$ws.Range("L103").Value = 2247
$ws.Range("M103").Value = 0.647
$ws.Range("N103").Value = 0.006
$ws.Range("O103").Value = 157.29

# Row 105
$ws.Range("C105").Value = 44063
$ws.Range("G105").Value = 105
$ws.Range("H105").Value = 87885
$ws.Range("I105").Value = 5.913
$ws.Range("J105").Value = 1256
$ws.Range("K105").Value = 0.085
$ws.Range("L105").Value = 1307
$ws.Range("M105").Value = 0.088
$ws.Range("N105").Value = 0.082
$ws.Range("O105").Value = 12.199

